$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.398115754127502
$ws.Range("B1").Value = 2.523339986801147
$ws.Range("C1").Value = 4.466956615447998
$ws.Range("D1").Value = 1.722910404205322
$ws.Range("E1").Value = 1.059664607048035
